$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 5 (Module Bluetooth / Iduino BT-Board HC-05), shifting rows 5-7 down to 6-8
$ws.Rows("5:5").Insert()

# Fill in the new row 5
$ws.Range("A5").Value = "Module Bluetooth"
$ws.Range("B5").Value = "Iduino BT-Board HC-05"
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = "https://www.aranacorp.com/fr/votre-arduino-communique-avec-le-module-hc-05/ "

# Replace the old rich-text "gotronic" cell (now row 8) with a plain-text value that has a trailing space
$ws.Range("D8").Value = "https://www.gotronic.fr/pj-496.pdf "

# Rebuild all hyperlinks (the row insert does not shift existing hyperlink anchors)
$ws.Cells.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("D2"), "https://store.arduino.cc/products/arduino-due", [Type]::Missing, [Type]::Missing, "https://store.arduino.cc/products/arduino-due ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "https://datasheetspdf.com/pdf-file/1251016/ETC/DYP-ME007/1", [Type]::Missing, [Type]::Missing, "https://datasheetspdf.com/pdf-file/1251016/ETC/DYP-ME007/1 ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D4"), "http://wiki.sunfounder.cc/index.php?title=Bluetooth_4.0_HM-10_Master_Slave_Module", [Type]::Missing, [Type]::Missing, "http://wiki.sunfounder.cc/index.php?title=Bluetooth_4.0_HM-10_Master_Slave_Module ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D5"), "https://www.aranacorp.com/fr/votre-arduino-communique-avec-le-module-hc-05/", [Type]::Missing, [Type]::Missing, "https://www.aranacorp.com/fr/votre-arduino-communique-avec-le-module-hc-05/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D6"), "http://www.ee.ic.ac.uk/pcheung/teaching/DE1_EE/stores/sg90_datasheet.pdf", [Type]::Missing, [Type]::Missing, "http://www.ee.ic.ac.uk/pcheung/teaching/DE1_EE/stores/sg90_datasheet.pdf ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D7"), "https://www.robotshop.com/content/ZIP/documentation-sabertooth-2x12-rc.zip", [Type]::Missing, [Type]::Missing, "https://www.robotshop.com/content/ZIP/documentation-sabertooth-2x12-rc.zip ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D8"), "https://www.gotronic.fr/pj-496.pdf", [Type]::Missing, [Type]::Missing, "https://www.gotronic.fr/pj-496.pdf ") | Out-Null

# Restore the active selection as in the target file
$ws.Range("B6").Select()
